$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.932.50'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.649.04'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.51%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.73'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3919'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3876'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '51.31'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.365'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.006'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08484'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.03'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.222'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.899'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +5.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001315'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.646.61'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '95.00'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06992'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.08'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.946'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.006'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.61%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.63'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.924.15'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'LidoDAOToken'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.122'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +8.89%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.456'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.19'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.50'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '139.57'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.317'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.877'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -6.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.499'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.831.09'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.031'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.50%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.76%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.08124'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.700'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '10.98'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +5.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2710'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.09153'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.81%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.7552'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.50'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.426'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6951'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.486'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.096'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08286'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '134.62'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.415'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +8.01%  '
